# Modified DSL for EB
# Clear the "Pass" values out of the Results column (J2:J5) on the
# TestCases sheet, and select that range afterwards (matches the
# author's saved selection in the edited workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$range = $ws.Range("J2:J5")
$range.ClearContents()

# Column J (Results) shrinks its cached "best fit" width now that the
# short "Pass" values are gone (only the "Results" header remains).
$ws.Columns.Item(10).ColumnWidth = 5.17

$ws.Activate()
$range.Select()
